$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 10 (Event column, A10): "When nav item is clicked" -> "When header is scrolled above section#about"
$ws.Range("A10").Value = "When header is scrolled above section#about"

# Update row 11 (Task column, B11): "Smooth scroll to section" -> "Re-attach nav menu in section#landing"
$ws.Range("B11").Value = "Re-attach nav menu in section#landing"

# Update row 12 (Event column, A12): "When header is scrolled above section#about" -> "Touch screen anchor click"
$ws.Range("A12").Value = "Touch screen anchor click"

# Update row 13 (Task column, B13): "Re-attach nav menu in section#landing" -> "Change font color back after click"
$ws.Range("B13").Value = "Change font color back after click"

# Remove the now-obsolete rows 14-17 entirely
$ws.Range("A14:B17").EntireRow.Delete() | Out-Null

# Update the selected cell to match the final workbook state
$ws.Range("A17").Select() | Out-Null
